$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 126698504
$ws.Range("E2").Value = 1633265
$ws.Range("G2").Value = 122584810
$ws.Range("I2").Value = 2480429
$ws.Range("J2").Value = 6212

$ws.Range("C3").Value = 217425850
$ws.Range("D3").Value = 17533
$ws.Range("E3").Value = 1546499
$ws.Range("F3").Value = 3
$ws.Range("G3").Value = 201115588
$ws.Range("H3").Value = 21652
$ws.Range("I3").Value = 14763763
$ws.Range("J3").Value = 15404

$ws.Range("C4").Value = 249360306
$ws.Range("D4").ClearContents()
$ws.Range("E4").Value = 2061014
$ws.Range("F4").ClearContents()
$ws.Range("G4").Value = 245331922
$ws.Range("H4").Value = 4381
$ws.Range("I4").Value = 1967370
$ws.Range("J4").Value = 5792

$ws.Range("C5").Value = 68711939
$ws.Range("E5").Value = 1355891
$ws.Range("G5").Value = 66471658
$ws.Range("H5").Value = 246
$ws.Range("I5").Value = 884390
$ws.Range("J5").Value = 10123

$ws.Range("C6").Value = 14365172
$ws.Range("D6").ClearContents()
$ws.Range("E6").Value = 28077
$ws.Range("G6").Value = 14188517
$ws.Range("I6").Value = 148578
$ws.Range("J6").Value = 99

$ws.Range("C7").Value = 12823515
$ws.Range("G7").Value = 12085476
$ws.Range("H7").Value = 4
$ws.Range("I7").Value = 479274
